$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# Row 20 already has "User Soft Deletes" in A20; mark it Done.
$ws.Range("B20").Value = "Done"

# Add a new row 21 for the new task.
$ws.Range("A21").Value = "Payment Transaction in Create Challenge"
$ws.Range("B21").Value = "Done"

$ws.Range("A14").Select()
$ws.Range("B20:B21").Select()
